$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are numeric-looking text (e.g. "1.010", "192.30") that
# Excel auto-converts to actual numbers on assignment, which would lose the exact
# textual formatting used by this sheet (inline/shared strings). Force them to stay
# text by applying a text number format before the write, then resetting the style
# back to Normal so no stray formatting is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.283.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.663.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5335"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("E8").Value = "  +1.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06359"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.672.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.891.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5536"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8189"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.280.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.665"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1226"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.236"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.578"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.616"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.823"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9591"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.431"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5814"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01607"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.878"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8540"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.94%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.009"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.047.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.804.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4375"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.942"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("E51").Value = "  +0.27%  "
